$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $style = $ws.Range($cellRef).Style
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = $style
}

Set-TextValue $ws "D2" "96.998.73"
Set-TextValue $ws "E2" "  -0.49%  "
Set-TextValue $ws "D3" "3.693.49"
Set-TextValue $ws "E3" "  +2.77%  "
Set-TextValue $ws "E4" "  +0.00%  "
Set-TextValue $ws "D5" "239.91"
Set-TextValue $ws "E5" "  -0.99%  "
Set-TextValue $ws "D6" "1.90"
Set-TextValue $ws "E6" "  +8.74%  "
Set-TextValue $ws "D7" "653.97"
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "D8" "0.427"
Set-TextValue $ws "E8" "  -1.60%  "
Set-TextValue $ws "E9" "  +3.03%  "
Set-TextValue $ws "D10" "1.00"
Set-TextValue $ws "E10" "  +0.05%  "
Set-TextValue $ws "D11" "3.692.51"
Set-TextValue $ws "E11" "  +2.76%  "
Set-TextValue $ws "D12" "45.51"
Set-TextValue $ws "E12" "  +2.07%  "
Set-TextValue $ws "D13" "0.207"
Set-TextValue $ws "E13" "  +0.98%  "
Set-TextValue $ws "D14" "6.91"
Set-TextValue $ws "E14" "  +6.50%  "
Set-TextValue $ws "D15" "4.378.78"
Set-TextValue $ws "E15" "  +2.76%  "
Set-TextValue $ws "E16" "  +2.31%  "
Set-TextValue $ws "D17" "96.738.59"
Set-TextValue $ws "E17" "  -0.41%  "
Set-TextValue $ws "D18" "9.08"
Set-TextValue $ws "E18" "  +4.05%  "
Set-TextValue $ws "D19" "3.678.67"
Set-TextValue $ws "E19" "  +2.38%  "
Set-TextValue $ws "D20" "19.36"
Set-TextValue $ws "E20" "  +6.31%  "
Set-TextValue $ws "D21" "12.87"
Set-TextValue $ws "E21" "  +1.95%  "
Set-TextValue $ws "E22" "  -0.66%  "
Set-TextValue $ws "D23" "526.91"
Set-TextValue $ws "E23" "  +1.58%  "
Set-TextValue $ws "D24" "3.51"
Set-TextValue $ws "E24" "  +0.09%  "
Set-TextValue $ws "D25" "7.15"
Set-TextValue $ws "E25" "  +2.75%  "
Set-TextValue $ws "E26" "  -2.11%  "
Set-TextValue $ws "D27" "102.13"
Set-TextValue $ws "E27" "  +0.02%  "
Set-TextValue $ws "D28" "13.45"
Set-TextValue $ws "E28" "  +2.10%  "
Set-TextValue $ws "E29" "  -1.31%  "
Set-TextValue $ws "D30" "12.60"
Set-TextValue $ws "E30" "  +4.07%  "
Set-TextValue $ws "E31" "  +1.80%  "
Set-TextValue $ws "E32" "  -0.17%  "
Set-TextValue $ws "E33" "  +13.96%  "
Set-TextValue $ws "E34" "  +0.50%  "
Set-TextValue $ws "D35" "32.77"
Set-TextValue $ws "E35" "  +2.30%  "
Set-TextValue $ws "E36" "  +0.10%  "
Set-TextValue $ws "D37" "0.615"
Set-TextValue $ws "E37" "  +7.04%  "
Set-TextValue $ws "D38" "656.20"
Set-TextValue $ws "E38" "  +6.12%  "
Set-TextValue $ws "D39" "9.08"
Set-TextValue $ws "E39" "  +3.58%  "
Set-TextValue $ws "D40" "7.07"
Set-TextValue $ws "E40" "  +16.82%  "
Set-TextValue $ws "E41" "  +5.60%  "
Set-TextValue $ws "D42" "2.03"
Set-TextValue $ws "E42" "  +4.14%  "
Set-TextValue $ws "D43" "0.968"
Set-TextValue $ws "E43" "  +3.83%  "
Set-TextValue $ws "D44" "38.62"
Set-TextValue $ws "E44" "  +16.86%  "
Set-TextValue $ws "D46" "0.457"
Set-TextValue $ws "E46" "  +7.58%  "
Set-TextValue $ws "E48" "  +0.28%  "
Set-TextValue $ws "D49" "8.82"
Set-TextValue $ws "E49" "  +2.60%  "
Set-TextValue $ws "D51" "3.57"
Set-TextValue $ws "E51" "  +2.14%  "
